# Commit: "Create functions to create clip-specific pages. Add BealeStreet"
#
# This script updates the abbreviation labels in column A for several films,
# switching the short codes used for clip-specific page generation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("film_info")

$ws.Range("A2").Value  = "AkeelahBee"
$ws.Range("A5").Value  = "DaysSummer"
$ws.Range("A8").Value  = "GoodWill"
$ws.Range("A10").Value = "PursuitHappyness"
$ws.Range("A14").Value = "LittleMiss"

# Move/record the active selection at the end of the data, as in the saved file
$ws.Range("A23").Select()
